$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the claim/siniestro numbers in column E.
# Prefix with an apostrophe so Excel keeps these numeric-looking values
# stored as text (preserving the existing cell style/number format and,
# for E3, the trailing spaces) instead of converting them to numbers.
$ws.Range("E2").Value = "'1120194100404"
$ws.Range("E3").Value = "'1120170200928  "

# Update the last selected cell in the sheet view
$ws.Range("G5").Select()
